$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 1499.8889
$ws.Range("I51").Value = 1499
$ws.Range("K51").Value = 1499
$ws.Range("M51").Value = -1015
$ws.Range("H55").Value = 552.3125
$ws.Range("I55").Value = 261.91666
$ws.Range("J55").Value = 1423.5
$ws.Range("K55").Value = 261.91666
$ws.Range("L55").Value = 1423.5
$ws.Range("M55").Value = -47.91665999999998
$ws.Range("N55").Value = -1851.5
$ws.Range("H92").Value = 851.95654
$ws.Range("I92").Value = 877.3158
$ws.Range("J92").Value = 731.5
$ws.Range("K92").Value = 877.3158
$ws.Range("L92").Value = 731.5
$ws.Range("M92").Value = 370.6842
$ws.Range("N92").Value = -3227.5
$ws.Range("H96").Value = 2543.3333
$ws.Range("I96").Value = 420.66666
$ws.Range("K96").Value = 1261.99998
$ws.Range("M96").Value = 111.0000199999999
$ws.Range("H101").Value = 442.5
$ws.Range("J101").Value = 485
$ws.Range("L101").Value = 1455
$ws.Range("N101").Value = -4699
$ws.Range("H103").Value = 1752.1818
$ws.Range("I103").Value = 1425.3334
$ws.Range("J103").Value = 1874.75
$ws.Range("K103").Value = 4276.0002
$ws.Range("L103").Value = 5624.25
$ws.Range("M103").Value = -3690.0002
$ws.Range("N103").Value = -6796.25
$ws.Range("H111").Value = 767.1429000000001
$ws.Range("I111").Value = 767.1429000000001
$ws.Range("K111").Value = 2301.4287
$ws.Range("M111").Value = 765.5712999999996
$ws.Range("H115").Value = 2999.25
$ws.Range("I115").Value = 2999.25
$ws.Range("K115").Value = 8997.75
$ws.Range("M115").Value = -7430.75
$ws.Range("H138").Value = 4402.5854
$ws.Range("J138").Value = 4995.6562
$ws.Range("L138").Value = 14986.9686
$ws.Range("N138").Value = -25266.9686
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1977.4445
$ws.Range("I2").Value = 1971
$ws.Range("K2").Value = 1971
$ws.Range("M2").Value = -1858
$ws.Range("H32").Value = 8236.275
$ws.Range("I32").Value = 7832.4863
$ws.Range("K32").Value = 7832.4863
$ws.Range("M32").Value = -7545.4863
$ws.Range("H97").Value = 2690.3333
$ws.Range("I97").Value = 612.375
$ws.Range("J97").Value = 6846.25
$ws.Range("K97").Value = 612.375
$ws.Range("L97").Value = 6846.25
$ws.Range("M97").Value = -116.375
$ws.Range("N97").Value = -7838.25
$ws.Range("H102").Value = 1843.6471
$ws.Range("I102").Value = 1871.375
$ws.Range("K102").Value = 1871.375
$ws.Range("M102").Value = -249.375
$ws.Range("H116").Value = 1977.4445
$ws.Range("I116").Value = 1971
$ws.Range("K116").Value = 1971
$ws.Range("M116").Value = 323
$ws.Range("H122").Value = 1439.1875
$ws.Range("I122").Value = 1201.8
$ws.Range("K122").Value = 3605.4
$ws.Range("M122").Value = -1155.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1977.4445
$ws.Range("I3").Value = 1971
$ws.Range("K3").Value = 1971
$ws.Range("M3").Value = -1857
$ws.Range("H86").Value = 3537.3125
$ws.Range("J86").Value = 5103.875
$ws.Range("L86").Value = 5103.875
$ws.Range("N86").Value = -7349.875
$ws.Range("H89").Value = 3537.3125
$ws.Range("J89").Value = 5103.875
$ws.Range("L89").Value = 25519.375
$ws.Range("N89").Value = -36751.375
$ws.Range("H94").Value = 2460.111
$ws.Range("I94").Value = 2752.5334
$ws.Range("K94").Value = 2752.5334
$ws.Range("M94").Value = -2301.5334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 47711.5
$ws.Range("J60").Value = 47711.5
$ws.Range("L60").Value = 47711.5
$ws.Range("N60").Value = -48733.5
$ws.Range("H108").Value = 56000
$ws.Range("J108").Value = 56000
$ws.Range("L108").Value = 56000
$ws.Range("N108").Value = -63680
$ws.Range("H134").Value = 1595
$ws.Range("I134").Value = 997.6667
$ws.Range("K134").Value = 2993.0001
$ws.Range("M134").Value = -458.0001000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 5047.5
$ws.Range("I103").Value = 95
$ws.Range("J103").Value = 10000
$ws.Range("K103").Value = 285
$ws.Range("L103").Value = 30000
$ws.Range("M103").Value = 594
$ws.Range("N103").Value = -31758
$ws.Range("H129").Value = 2892.5
$ws.Range("I129").Value = 1997.5
$ws.Range("K129").Value = 5992.5
$ws.Range("M129").Value = -992.5
$ws.Range("H136").Value = 2335.5
$ws.Range("I136").Value = 1956
$ws.Range("K136").Value = 5868
$ws.Range("M136").Value = -768

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1348.125
$ws.Range("I97").Value = 1348.125
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1348.125
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -852.125
$ws.Range("N97").ClearContents()
$ws.Range("H136").Value = 30000
$ws.Range("J136").Value = 30000
$ws.Range("L136").Value = 90000
$ws.Range("N136").Value = -95100

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 837.8
$ws.Range("I100").Value = 837.8
$ws.Range("K100").Value = 837.8
$ws.Range("M100").Value = -296.8
$ws.Range("H136").Value = 3063
$ws.Range("I136").Value = 2137.625
$ws.Range("J136").Value = 4913.75
$ws.Range("K136").Value = 6412.875
$ws.Range("L136").Value = 14741.25
$ws.Range("M136").Value = -3862.875
$ws.Range("N136").Value = -19841.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H100").Value = 790.7646999999999
$ws.Range("I100").Value = 818.5
$ws.Range("J100").Value = 751.1429000000001
$ws.Range("K100").Value = 1637
$ws.Range("L100").Value = 1502.2858
$ws.Range("M100").Value = -1096
$ws.Range("N100").Value = -2584.2858
$ws.Range("H113").Value = 502.94116
$ws.Range("I113").Value = 509.33334
$ws.Range("K113").Value = 1528.00002
$ws.Range("M113").Value = 641.9999800000001
